$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14 / Row 15: mark as tardy (NO OF HOURS LATE -> 1) and highlight the rows red ---
$ws.Range("I14").Value = 1
$ws.Range("I15").Value = 1

# Apply the "red" row color from the new color palette (cyan / orange / red) that was
# added to the workbook's styles. Rows 14-15 use the red swatch (FFDF5E5E).
$ws.Range("A14:J15").Interior.Color = 6184671

# --- Row 19: B19 ("TOTAL LEAVES ACCUMULATED" row) switches from a blank text marker to a boolean FALSE ---
# B19 lives inside the merged range A19:G19, so it has to be unmerged momentarily in
# order to address it directly; only B19's content changes, everything else in the
# merged area keeps its original value.
$ws.Range("A19:G19").UnMerge()
$ws.Range("B19").Value = $False

# --- Simplify the FLOOR(...) calls: drop the redundant third argument ---
$ws.Range("B22").Formula = '=FLOOR(F17/8,1)&"."&FLOOR(MOD(F17,8),1)&"."&(MOD(F17,8)-FLOOR(MOD(F17,8),1))*60'
$ws.Range("B23").Formula = '=FLOOR(H19,1)&"."&(H19-FLOOR(H19,1))*8&".0"'
$ws.Range("B24").Formula = '=FLOOR(I19,1)&"."&(I19-FLOOR(I19,1))*8&".0"'
$ws.Range("B27").Formula = '=FLOOR(K27/8,1)&"."&FLOOR(MOD(K27,8),1)&"."&(MOD(K27,8)-FLOOR(MOD(K27,8),1))*60'
